# Actualizo archivos de brechas de ingresos por región
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Sexo" column values were relabeled: Masculino -> Hombre, Femenino -> Mujer.
# Cells.Replace (LookAt:=xlWhole, i.e. 1) rewrites the shared-string text
# in place, so every row that referenced that shared string (34 rows each)
# keeps pointing at the same entry instead of spawning new shared-string
# records - matching how the workbook actually came back from Excel.
[void]$ws.Cells.Replace("Masculino", "Hombre", 1)
[void]$ws.Cells.Replace("Femenino", "Mujer", 1)

# Move the live selection/active cell to T73 (it was A66:XFD67 / A66 before).
[void]$ws.Range("T73").Select()
